$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "33.878.57"
Set-TextValue "E2" "  -2.14%  "
Set-TextValue "D3" "1.780.26"
Set-TextValue "E3" "  -0.56%  "
Set-TextValue "E4" "  -0.21%  "
Set-TextValue "D5" "221.08"
Set-TextValue "E5" "  -1.43%  "
Set-TextValue "D6" "0.551"
Set-TextValue "E6" "  -1.58%  "
Set-TextValue "D7" "0.998"
Set-TextValue "E7" "  -0.13%  "
Set-TextValue "D8" "31.09"
Set-TextValue "E8" "  -4.55%  "
Set-TextValue "E9" "  +0.65%  "
Set-TextValue "D10" "0.0709"
Set-TextValue "E10" "  +6.10%  "
Set-TextValue "E11" "  -1.64%  "
Set-TextValue "D12" "2.034.77"
Set-TextValue "E12" "  -0.58%  "
Set-TextValue "D13" "1.775.88"
Set-TextValue "E13" "  -0.63%  "
Set-TextValue "E14" "  -4.74%  "
Set-TextValue "E15" "  -1.20%  "
Set-TextValue "D16" "33.863.89"
Set-TextValue "E16" "  -2.20%  "
Set-TextValue "D17" "4.22"
Set-TextValue "E17" "  -1.59%  "
Set-TextValue "D18" "68.01"
Set-TextValue "E18" "  -0.99%  "
Set-TextValue "D19" "244.81"
Set-TextValue "E19" "  -3.52%  "
Set-TextValue "D20" "0.0₃0776"
Set-TextValue "E20" "  +1.31%  "
Set-TextValue "E21" "  -0.07%  "
Set-TextValue "D22" "10.65"
Set-TextValue "E22" "  +2.37%  "
Set-TextValue "D23" "4.07"
Set-TextValue "E23" "  -3.73%  "
Set-TextValue "D24" "2.09"
Set-TextValue "E24" "  -1.98%  "
Set-TextValue "E25" "  -1.31%  "
Set-TextValue "D26" "16.39"
Set-TextValue "E26" "  +0.13%  "
Set-TextValue "E27" "  -1.28%  "
Set-TextValue "D29" "0.998"
Set-TextValue "E29" "  -0.30%  "
Set-TextValue "E30" "  +0.78%  "
Set-TextValue "E31" "  -1.48%  "
Set-TextValue "E32" "  +0.26%  "
Set-TextValue "D33" "3.50"
Set-TextValue "E33" "  -2.24%  "
Set-TextValue "D34" "1.81"
Set-TextValue "E34" "  -1.95%  "
Set-TextValue "D35" "1.398.26"
Set-TextValue "E35" "  -3.07%  "
Set-TextValue "E37" "  -0.40%  "
Set-TextValue "E38" "  -1.95%  "
Set-TextValue "D39" "0.934"
Set-TextValue "E39" "  +3.81%  "
Set-TextValue "B40" "Aave"
Set-TextValue "C40" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D40" "79.14"
Set-TextValue "E40" "  -4.71%  "
Set-TextValue "B41" "HuobiToken"
Set-TextValue "C41" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D41" "2.34"
Set-TextValue "E41" "  -0.37%  "
Set-TextValue "E42" "  -3.64%  "
Set-TextValue "E43" "  +1.69%  "
Set-TextValue "E44" "  +0.52%  "
Set-TextValue "E45" "  -2.92%  "
Set-TextValue "E46" "  -1.43%  "
Set-TextValue "D47" "1.933.90"
Set-TextValue "E47" "  -0.10%  "
Set-TextValue "D48" "104.90"
Set-TextValue "E48" "  +1.64%  "
Set-TextValue "D49" "0.995"
Set-TextValue "E49" "  -0.48%  "
Set-TextValue "D50" "11.71"
Set-TextValue "E50" "  -2.55%  "
Set-TextValue "E51" "  -2.03%  "
